# Updated analysed sessions Nathan
# Bumps the "Output_id" (column B) count for the first 8 "surround-mod" sessions
# and flags them with a new "Comments" column (P) noting "dt_post_stim is 0".
# Also leaves a quote-prefixed (text-formatted) blank cell at A23, matching the
# author's scratch annotation below the table, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Output_id) increments for rows 2-9 ---
$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 9
$ws.Range("B9").Value = 5

# --- New "Comments" column (P) ---
$ws.Range("P1").Value = "Comments"
$ws.Range("P2").Value = "dt_post_stim is 0"
$ws.Range("P3").Value = "dt_post_stim is 0"
$ws.Range("P4").Value = "dt_post_stim is 0"
$ws.Range("P5").Value = "dt_post_stim is 0"
$ws.Range("P6").Value = "dt_post_stim is 0"
$ws.Range("P7").Value = "dt_post_stim is 0"
$ws.Range("P8").Value = "dt_post_stim is 0"
$ws.Range("P9").Value = "dt_post_stim is 0"

# --- Blank, quote-prefixed (left-aligned text) cell left at A23 ---
$ws.Range("A23").Value = "'"
$ws.Range("A23").Value = ""

# --- Restore the author's last on-screen selection ---
$ws.Range("E28").Select()
